$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = "['MEC-3A-M.S.R. ar Cond.', -, -, -]"
$ws.Range("D11").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("D14").Value = "[Tiago Freitas-M.S.R. ar Cond.-3A, -, -, -]"
$ws.Range("D15").Value = "['MEC-3A-M.S.R. ar Cond.', -, -, -]"
$ws.Range("F16").Value = "[-, -, -, 'MEC-3A-M.S.R. ar Cond.']"
